# V2G_ratio.xlsx edit: rotate the scenario labeling by one step
# (No V2G -> Low -> Moderate -> High -> V2G mandate -> No V2G) and
# correspondingly rotate the BEV "value" data so each scenario keeps the
# data that actually belongs to it after the relabeling.
#
# "Early" is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rotate the Scenario labels (column B) for every Drive_train block.
#    Old block row ranges (1950-2050, 4 drive trains x 101 rows each):
#      Low          rows 2-405
#      Moderate     rows 406-809
#      High         rows 810-1213
#      V2G mandate  rows 1214-1617
#      No V2G       rows 1618-2021
#      Early        rows 2022-2425 (unchanged)
# ---------------------------------------------------------------------

$ws.Range("B2:B405").Value     = "No V2G"
$ws.Range("B406:B809").Value   = "Low"
$ws.Range("B810:B1213").Value  = "Moderate"
$ws.Range("B1214:B1617").Value = "High"
$ws.Range("B1618:B2021").Value = "V2G mandate"

# ---------------------------------------------------------------------
# 2) Rotate the "value" column (E) for the BEV sub-blocks only, so the
#    data travels with its real scenario identity (the diff shows only
#    BEV rows changing -- ICE/PHEV/H values are all 0 in every scenario
#    so no visible difference, but we still only touch BEV to match).
#
#    Old BEV block rows:
#      Low BEV          103-203
#      Moderate BEV      507-607
#      High BEV           911-1011
#      V2G mandate BEV   1315-1415
#      No V2G BEV        1719-1819
#
#    Snapshot all five blocks first (read-only), then write them back in
#    rotated order so none of the source data is clobbered mid-flight.
# ---------------------------------------------------------------------

$lowBev        = $ws.Range("E103:E203").Value2
$moderateBev   = $ws.Range("E507:E607").Value2
$highBev       = $ws.Range("E911:E1011").Value2
$v2gMandateBev = $ws.Range("E1315:E1415").Value2
$noV2gBev      = $ws.Range("E1719:E1819").Value2

$ws.Range("E103:E203").Value    = $noV2gBev       # now labeled "No V2G"  -> gets the real No V2G data
$ws.Range("E507:E607").Value    = $lowBev         # now labeled "Low"     -> gets the real Low data
$ws.Range("E911:E1011").Value   = $moderateBev    # now labeled "Moderate"-> gets the real Moderate data
$ws.Range("E1315:E1415").Value  = $highBev        # now labeled "High"    -> gets the real High data
$ws.Range("E1719:E1819").Value  = $v2gMandateBev  # now labeled "V2G mandate" -> gets the real V2G mandate data
